# This script swaps the values of columns B and D:G between pairs of
# adjacent rows in the active worksheet. Columns A (serial no.) and C
# (item name) stay put; B (item code), D (rate), E (mrp), F (qty) and
# G (value) are exchanged between the two rows of each pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    313,314,
    355,356,
    372,373,
    400,401,
    421,422,
    457,458,
    579,580,
    583,584,
    586,587,
    593,594,
    599,600,
    604,605,
    687,688,
    709,710,
    720,721,
    859,860,
    889,890
)

for ($i = 0; $i -lt $rowPairs.Length; $i += 2) {
    $r1 = $rowPairs[$i]
    $r2 = $rowPairs[$i + 1]

    $range1 = $ws.Range("B$r1" + ":G$r1")
    $range2 = $ws.Range("B$r2" + ":G$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
